$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (inventory counts) per the revenue fix for manager class
$ws.Range("A2").Value = 996
$ws.Range("B2").Value = 959
$ws.Range("C2").Value = 959
$ws.Range("D2").Value = 959
$ws.Range("E2").Value = 966
$ws.Range("F2").Value = 966
$ws.Range("G2").Value = 996
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
